$d = $word.ActiveDocument

# The "Solutions" table is the only table in the document; grab it and
# append a brand-new row after the existing "Problem 6" row, mirroring
# the row structure used throughout the table (Problem / Part / Solution).
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "7"
$newRow.Cells.Item(2).Range.Text = "-"
$newRow.Cells.Item(3).Range.Text = "0.377"
